$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / add Activities text (column C) --------------------------
# Order matters: it determines how the shared-strings table is laid out.
# 1) "Meeting" is a brand-new activity typed into a previously blank cell,
#    so it gets appended as a new shared string first.
$ws.Range("C4").Value = "Meeting"

# 2) The first two existing activity rows get their text revised in place
#    (they keep their original shared-string slots, just with new text).
$ws.Range("C2").Value = "Design - Setting up Github account, creating templates for progress reports, creating project timeline system."
$ws.Range("C3").Value = "Design - Revising project timeline system. "

# 3) The remaining new activity rows are filled in.
$ws.Range("C5").Value = "Design - Perfoming Requirement Analysis"
$ws.Range("C6").Value = "Design - Revising, writing Introduction, and constructing ER Diagrams"
$ws.Range("C7").Value = "Design - Revising"
$ws.Range("C8").Value = "Design - Revising ER Diagrams"
$ws.Range("C9").Value = "Design - Creating cover page and revising"

# --- Fill in Date (column A) and Hours (column B) for the new rows ----
$ws.Range("A4").Value = 42755
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = 42760
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = 42765
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = 42767
$ws.Range("B7").Value = 1

$ws.Range("A8").Value = 42770
$ws.Range("B8").Value = 1

$ws.Range("A9").Value = 42771
$ws.Range("B9").Value = 1

# Give the newly-populated date cells (A4:A9) the same date number format
# already used by A3, without disturbing the values we just set.
$ws.Range("A3").Copy()
$ws.Range("A4:A9").PasteSpecial(-4122) # xlPasteFormats

# --- Column C is now narrower (the activity text is shorter) ----------
$ws.Columns.Item(3).ColumnWidth = 53.5

# --- The active selection moved to C9 ----------------------------------
$ws.Range("C9").Select()
